{"js": "// Replace the date line and the 25 two-digit multiplication problems\n// with their new values, as described by the diff.\nconst replacements = [\n  [\"2024-03-04 Monday\", \"2024-03-05 Tuesday\"],\n  [\"47\\u00d714=\", \"59\\u00d753=\"],\n  [\"42\\u00d727=\", \"80\\u00d744=\"],\n  [\"45\\u00d796=\", \"59\\u00d712=\"],\n  [\"38\\u00d781=\", \"67\\u00d744=\"],\n  [\"68\\u00d772=\", \"27\\u00d724=\"],\n  [\"91\\u00d711=\", \"29\\u00d761=\"],\n  [\"72\\u00d746=\", \"84\\u00d771=\"],\n  [\"42\\u00d770=\", \"17\\u00d718=\"],\n  [\"94\\u00d776=\", \"48\\u00d780=\"],\n  [\"85\\u00d740=\", \"99\\u00d792=\"],\n  [\"90\\u00d753=\", \"98\\u00d731=\"],\n  [\"59\\u00d764=\", \"51\\u00d752=\"],\n  [\"24\\u00d727=\", \"11\\u00d782=\"],\n  [\"31\\u00d784=\", \"50\\u00d731=\"],\n  [\"80\\u00d723=\", \"59\\u00d756=\"],\n  [\"97\\u00d794=\", \"87\\u00d715=\"],\n  [\"99\\u00d793=\", \"47\\u00d772=\"],\n  [\"29\\u00d790=\", \"84\\u00d745=\"],\n  [\"91\\u00d766=\", \"58\\u00d718=\"],\n  [\"82\\u00d772=\", \"92\\u00d771=\"],\n  [\"40\\u00d757=\", \"60\\u00d767=\"],\n  [\"44\\u00d779=\", \"54\\u00d718=\"],\n  [\"67\\u00d789=\", \"76\\u00d719=\"],\n  [\"57\\u00d755=\", \"18\\u00d738=\"],\n  [\"94\\u00d759=\", \"77\\u00d788=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 two-digit multiplication problems\n# with their new values, as described by the diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-04 Monday\", \"2024-03-05 Tuesday\"),\n    @(\"47\u00d714=\", \"59\u00d753=\"),\n    @(\"42\u00d727=\", \"80\u00d744=\"),\n    @(\"45\u00d796=\", \"59\u00d712=\"),\n    @(\"38\u00d781=\", \"67\u00d744=\"),\n    @(\"68\u00d772=\", \"27\u00d724=\"),\n    @(\"91\u00d711=\", \"29\u00d761=\"),\n    @(\"72\u00d746=\", \"84\u00d771=\"),\n    @(\"42\u00d770=\", \"17\u00d718=\"),\n    @(\"94\u00d776=\", \"48\u00d780=\"),\n    @(\"85\u00d740=\", \"99\u00d792=\"),\n    @(\"90\u00d753=\", \"98\u00d731=\"),\n    @(\"59\u00d764=\", \"51\u00d752=\"),\n    @(\"24\u00d727=\", \"11\u00d782=\"),\n    @(\"31\u00d784=\", \"50\u00d731=\"),\n    @(\"80\u00d723=\", \"59\u00d756=\"),\n    @(\"97\u00d794=\", \"87\u00d715=\"),\n    @(\"99\u00d793=\", \"47\u00d772=\"),\n    @(\"29\u00d790=\", \"84\u00d745=\"),\n    @(\"91\u00d766=\", \"58\u00d718=\"),\n    @(\"82\u00d772=\", \"92\u00d771=\"),\n    @(\"40\u00d757=\", \"60\u00d767=\"),\n    @(\"44\u00d779=\", \"54\u00d718=\"),\n    @(\"67\u00d789=\", \"76\u00d719=\"),\n    @(\"57\u00d755=\", \"18\u00d738=\"),\n    @(\"94\u00d759=\", \"77\u00d788=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
